# Update the cryptos list worksheet with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # Force the cell to keep its value as plain text, matching the
    # original inline-string cells (avoids Excel auto-coercing
    # number-looking strings such as "57.05" into numeric values).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.279.40"
$ws.Range("E2").Value = "  -1.42%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.049.93"
$ws.Range("E3").Value = "  -1.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "229.53"
$ws.Range("E5").Value = "  -1.63%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -1.05%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - Solana
Set-TextValue "D8" "57.05"
$ws.Range("E8").Value = "  -3.17%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.383"
$ws.Range("E9").Value = "  -2.67%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -2.24%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.27%  "

# Row 12 - Chainlink
Set-TextValue "D12" "14.75"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13 - Wrapped liquid staked Ether 2.0
$ws.Range("D13").Value = "2.350.13"
$ws.Range("E13").Value = "  -1.52%  "

# Row 14 - Avalanche
Set-TextValue "D14" "20.49"
$ws.Range("E14").Value = "  -2.99%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.753"
$ws.Range("E15").Value = "  -2.91%  "

# Row 16 - Polkadot
Set-TextValue "D16" "5.25"
$ws.Range("E16").Value = "  -1.97%  "

# Row 17 - Wrapped Ether
$ws.Range("D17").Value = "2.048.34"
$ws.Range("E17").Value = "  -2.72%  "

# Row 18 - Wrapped BTC
$ws.Range("D18").Value = "37.242.93"
$ws.Range("E18").Value = "  -1.25%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  -3.00%  "

# Row 20 - Litecoin
Set-TextValue "D20" "69.57"
$ws.Range("E20").Value = "  -2.78%  "

# Row 21 - Shiba Inu
$ws.Range("D21").Value = "0.0₃0821"
$ws.Range("E21").Value = "  -1.90%  "

# Row 22 - Bitcoin Cash
Set-TextValue "D22" "225.88"
$ws.Range("E22").Value = "  -1.10%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.06%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.10%  "

# Row 25 - PancakeSwap
Set-TextValue "D25" "2.32"
$ws.Range("E25").Value = "  -3.71%  "

# Row 26 - Cosmos
Set-TextValue "D26" "9.47"
$ws.Range("E26").Value = "  +2.82%  "

# Row 27 - Monero
Set-TextValue "D27" "168.29"
$ws.Range("E27").Value = "  -1.43%  "

# Row 28 - Kaspa
Set-TextValue "D28" "0.129"
$ws.Range("E28").Value = "  -4.04%  "

# Row 29 - Ethereum Classic
Set-TextValue "D29" "19.11"
$ws.Range("E29").Value = "  -2.11%  "

# Row 30 - Immutable X
Set-TextValue "D30" "1.35"
$ws.Range("E30").Value = "  -5.46%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.122"
$ws.Range("E31").Value = "  +0.56%  "

# Row 32 - Filecoin
Set-TextValue "D32" "4.52"
$ws.Range("E32").Value = "  -3.82%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.0623"
$ws.Range("E33").Value = "  -2.03%  "

# Row 34 - Internet Computer (DFINITY)
Set-TextValue "D34" "4.56"
$ws.Range("E34").Value = "  -3.28%  "

# Row 35 - Lido DAO Token
Set-TextValue "D35" "2.47"
$ws.Range("E35").Value = "  -1.10%  "

# Row 36 - WEMIX Token
Set-TextValue "D36" "1.81"
$ws.Range("E36").Value = "  -0.43%  "

# Row 37 - Render Token
Set-TextValue "D37" "3.28"
$ws.Range("E37").Value = "  -3.60%  "

# Row 38 - Binance USD
Set-TextValue "D38" "0.999"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39 - THORChain
$ws.Range("E39").Value = "  -2.57%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +4.42%  "

# Row 41 & 42 - Maker and Aave swapped positions/ranking
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D41" "97.74"
$ws.Range("E41").Value = "  -2.03%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.485.47"
$ws.Range("E42").Value = "  +3.03%  "

# Row 43 - Huobi Token
Set-TextValue "D43" "2.89"
$ws.Range("E43").Value = "  +0.37%  "

# Row 44 - Cronos
Set-TextValue "D44" "0.0950"
$ws.Range("E44").Value = "  -2.93%  "

# Row 45 - Trust Wallet Token
$ws.Range("E45").Value = "  +2.23%  "

# Row 46 - Injective Protocol
Set-TextValue "D46" "16.70"
$ws.Range("E46").Value = "  -0.16%  "

# Row 47 - FTX Token
$ws.Range("E47").Value = "  -3.13%  "

# Row 48 - ARBITRUM
Set-TextValue "D48" "1.03"

# Row 49 - Frax Share
Set-TextValue "D49" "7.19"
$ws.Range("E49").Value = "  -3.31%  "

# Row 50 - MX Token
$ws.Range("E50").Value = "  -1.98%  "

# Row 51 - Rocket Pool ETH
$ws.Range("D51").Value = "2.237.05"
$ws.Range("E51").Value = "  -1.49%  "
